$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add missing Edm Type Decimal: a new "decimal[decimal]" column (J) with
# sample decimal values for the two data rows.
$ws.Range("J1").Value = "decimal[decimal]"
$ws.Range("J2").Value = 15.3
$ws.Range("J3").Value = 14.3

# Explicitly (re-)apply the General number format on the last value so the
# cell carries its own style record (matches the other typed columns).
$ws.Range("J3").NumberFormat = "General"

# Move the active selection to the newly filled-in cell.
$ws.Range("J3").Select()
